$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated supplier summary data (rows 2-11), reflecting latest P&L calculation
$data = @(
    @("S001","fuel",39482.46000000001,4,1974.123000000001,"fine"),
    @("S009","catering",92287.27,1,8389.751818181818,"fine"),
    @("S006","maintenance",231730.48,2,11034.78476190476,"investigate"),
    @("S005","maintenance",277594.74,1,9914.097857142857,"fine"),
    @("S007","catering",126003.02,0,7411.942352941175,"fine"),
    @("S004","maintenance",267678.0900000001,3,14871.005,"investigate"),
    @("S002","fuel",48666.21000000002,0,2027.758750000001,"fine"),
    @("S003","fuel",46906.63000000002,2,2039.418695652175,"fine"),
    @("S008","catering",157532.35,1,8291.176315789473,"fine"),
    @("S010","catering",174173.66,4,8708.683000000001,"fine")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $rowIndex++
}
